# Apply attendance computation updates to the active worksheet.
# Sets specific cells in columns D, E, G, H from 0 to 1 as per the
# commit "final code with comment / tut 6".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: mark as Invalid and Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-6: mark as Real attendance (Total Attendance Count + Real)
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Rows 7-11: mark as Absent
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# Row 12: mark as Real attendance
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Rows 13-18: mark as Absent
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
